# Update Name of Algo
# Apply updated imputed values to the RandomForest result data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new value
$updates = @{
    "E7"   = 15.41209999999999
    "D10"  = -9.11399999999999
    "D12"  = -7.130299999999997
    "E15"  = 15.8363
    "D18"  = -8.059199999999993
    "E20"  = 15.8107
    "E29"  = 17.12360000000001
    "E30"  = 15.3993
    "E31"  = 15.78730000000001
    "D37"  = -7.579199999999995
    "E40"  = 17.01920000000002
    "D55"  = -8.289900000000001
    "D68"  = -7.4585
    "E68"  = 17.15170000000002
    "E76"  = 16.28519999999999
    "D77"  = -5.680599999999999
    "D78"  = -7.521700000000005
    "E87"  = 16.3718
    "E88"  = 16.4028
    "E96"  = 16.22949999999998
    "E98"  = 15.3858
    "E101" = 16.74160000000002
    "E102" = 16.73379999999999
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
